$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old single-column data -----------------------------------
$ws.Cells.Clear()

# --- New table data: Priority | Name | Status ---------------------------
$data = @(
    @("Priority", "Name", "Status"),
    @(2, "Report cards have default behavior", "pending"),
    @(1, "Refactor css classes for link cards", "done"),
    @(6, "Buttons have images", "pending"),
    @(3, "Moving an event is reflected in the reports", "pending"),
    @(4, "Report fields are appealing", "pending"),
    @(5, "The 'today' button works correctly", "pending"),
    @(7, "Calendar navigation buttons match the card style", "pending"),
    @(1, "Refactor controller (architecture)", "current"),
    @(3, "Report cards shows data for whole day, when calendar view is day", "pending")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- Column B now holds the wide text, matching the old column A style ---
$ws.Columns.Item(1).ColumnWidth = 8.43
$ws.Columns.Item(2).ColumnWidth = 67.5703125

# --- Apply an AutoFilter over the table, filtering Status -----------------
$rng = $ws.Range("A1:C9")
$rng.AutoFilter(3, "pending")

# --- Selection, matching the post-edit state ------------------------------
$ws.Range("C11").Select()

# --- Page setup: portrait orientation --------------------------------------
$ws.PageSetup.Orientation = 1

$wb.Save()
